$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor name entered on the form (next to "Supervisor Name:")
$ws.Range("G6").Value = "Ankita Gangotra"

# Supervisor sign-off: initials + date
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").Value = 41800
$ws.Range("D27").NumberFormat = "m/d/yy"

# Leave the selection where the user last clicked (the supervisor's date cell)
$ws.Range("D27:E27").Select()
